$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.07891583442688
$ws.Range("B1").Value = 2.119692325592041
$ws.Range("C1").Value = 9.061097145080566
$ws.Range("D1").Value = 1.033445239067078
$ws.Range("E1").Value = 0.9784101247787476
